# agregar tiempo a la ventana checklist
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the ticket/time number shown in the checklist window (cell W2).
# Leading apostrophe forces the numeric-looking value to be stored as text
# (matching the existing text-typed "Numero Propuesta" cell), then the
# style is reset back to Normal/default so no extra formatting sticks.
$ws.Range("W2").Value = "'4900028"
$ws.Range("W2").Style = "Normal"
